$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 56: new task entry "BA: SBERT (Hf)" dated 45203, keep existing D56 content untouched
# (set before row 55 so the shared-string table registers "BA: SBERT (Hf)" ahead of "BA: USE (google)")
$ws.Cells.Item(56, 1).Value = 45203
$ws.Cells.Item(56, 2).Value = "BA: SBERT (Hf)"

# Row 55: new task entry "BA: USE (google)" dated 45202, keep existing D55 content untouched
$ws.Cells.Item(55, 1).Value = 45202
$ws.Cells.Item(55, 2).Value = "BA: USE (google)"

# Copy style (number format, borders, alignment) from row 54's A/B cells onto the new A55/B55 and A56/B56 cells
$ws.Cells.Item(54, 1).Copy()
$ws.Cells.Item(55, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(56, 1).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(54, 2).Copy()
$ws.Cells.Item(55, 2).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(56, 2).PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Remove old rows 61 and 62 (leftover "infersent" TODO notes that were dropped)
$ws.Rows.Item(62).Delete() | Out-Null
$ws.Rows.Item(61).Delete() | Out-Null

# Update the view to reflect where the author was last working
$excel.Goto($ws.Range("A47"), $true) | Out-Null
$ws.Range("D57").Select() | Out-Null

$wb.Save()
